$wb = $excel.ActiveWorkbook

# --- Sheet1: performance_params_0 ---
$ws1 = $wb.Worksheets.Item("performance_params_0")
$ws1.Range("C2").Value = 2836
$ws1.Range("D2").Value = 2836

# --- Sheet2: Scaling ---
$ws2 = $wb.Worksheets.Item("Scaling")

# Add new "Optimal" / "Threshold" header columns, copying the header style
# used by the existing D1 header cell.
$ws2.Range("D1").Copy()
$ws2.Range("E1:F1").PasteSpecial(-4122)
$ws2.Range("E1").Value = "Optimal"
$ws2.Range("F1").Value = "Threshold"
$excel.CutCopyMode = 0

$ws2.Range("B2").Value = 2316
$ws2.Range("C2").Value = 2971

$ws2.Range("B3").Value = 1.7
$ws2.Range("C3").Value = 6

$ws2.Range("B4").Value = 45.7
$ws2.Range("C4").Value = 53

# --- Selection / active sheet & tab ---
$ws1.Range("D9").Select()
$ws1.Activate()
